$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column I: the "=" note was replaced by real numeric data (hour/min/sec
#     split of the new clock duration calculation). ---
$ws.Range("I3").Value = 0
$ws.Range("I4").Value = 1
$ws.Range("I5").Value = 2939

# I6 mirrors G6/H6: a time-of-day fraction formatted as h:mm.
$ws.Range("I6").NumberFormat = "h:mm"
$ws.Range("I6").Value = 0.041666666666666664

# --- Row 4: total minutes (seconds -> combine hh/mm/ss raw columns) ---
$ws.Range("J4").Formula = "=(G4+G5)"
$ws.Range("K4").Formula = "=(H4+H5)"
$ws.Range("L4").Formula = "=(I4+I5)"

# --- Row 6: convert each of those totals back down to the next unit (/60),
#     displayed with one decimal place. ---
$ws.Range("J6:L6").NumberFormat = "0.0"
$ws.Range("J6").Formula = "=J4/60"
$ws.Range("K6").Formula = "=K4/60"
$ws.Range("L6").Formula = "=L4/60"

# --- Row 3: add the row-6 carry back onto the original hour/min/sec,
#     displayed as whole numbers. ---
$ws.Range("J3:L3").NumberFormat = "0"
$ws.Range("J3").Formula = "=G3+J6"
$ws.Range("K3").Formula = "=H3+K6"
$ws.Range("L3").Formula = "=I3+L6"

# Column J widened to fit the new numbers (closest attainable width to the
# author's manual resize).
$ws.Columns(10).ColumnWidth = 11.6

# Author finished by leaving the selection on L4.
$ws.Range("L4").Select() | Out-Null
